$d = $word.ActiveDocument
$ins = $d.Range(0, 0)
$ins.InsertBefore("X")
Write-Host "after insert, text[0..3]=" $d.Range(0,3).Text
$rng = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $rng)
# Now delete the temp char X at position 0
$d.Range(0,1).Delete()
Write-Host "final text[0..10]=" $d.Range(0,10).Text
